$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 9) below the existing data (rows 1-8).
$ws.Range("A9").Value = 9597.3799999999992
$ws.Range("B9").Value = 9794.24
$ws.Range("C9").Value = 79.650000000000006
$ws.Range("D9").Value = 78.05
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -2.0099999999999998
$ws.Range("G9").Value = 42612.672905092593
$ws.Range("H9").Value = $false

# Match the date/time number format already used in column G (e.g. G8) -
# this is Excel's built-in format #22 ("m/d/yy h:mm"), so re-applying the
# same format string reuses the existing style instead of creating a new one.
$ws.Range("G9").NumberFormat = "m/d/yy h:mm"
